$p = $ppt.ActivePresentation

# --- 1. Table on slide 16 switches from the custom "Table_0" style to the
#        built-in "Medium Style 2 - Accent 3" table style. ---
$s = $p.Slides.Item(16)
$shp = $s.Shapes.Item(3)
$tbl = $shp.Table
$tbl.ApplyStyle("{5928D527-6A76-4F2C-89E7-7EE7A0985117}")

# --- 2. Re-colour the deck's theme from the custom "Integral" palette to the
#        stock Office palette (dk1/lt1 are already black/white in both). ---
$sm = $p.SlideMaster
$cs = $sm.Theme.ThemeColorScheme
# dk2 -> 44546A
$cs.Item(3).RGB  = 6968388
# lt2 -> E7E6E6
$cs.Item(4).RGB  = 15132391
# accent1 -> 5B9BD5
$cs.Item(5).RGB  = 13998939
# accent2 -> ED7D31
$cs.Item(6).RGB  = 3243501
# accent3 -> A5A5A5
$cs.Item(7).RGB  = 10855845
# accent4 -> FFC000
$cs.Item(8).RGB  = 49407
# accent5 -> 4472C4
$cs.Item(9).RGB  = 12874308
# accent6 -> 70AD47
$cs.Item(10).RGB = 4697456
# hlink -> 0563C1
$cs.Item(11).RGB = 12673797
# folHlink -> 954F72
$cs.Item(12).RGB = 7491477
